# SG_UBA1_VAR.docx template update ("Release 1.0")
#
# The underlying visible text of the long "DNA is analysed..." /
# "Somatic variant categorisation..." paragraph is unchanged by the
# commit (only the internal run-splitting / spell-check wrapper markup
# was reshuffled by Word on save), so no Find/Replace is required for
# that block. The two genuine content edits are:
#   1. The CV% statistics table of VAF bands/averages.
#   2. The cached SAVEDATE field result ("Reported" date).

$d = $word.ActiveDocument

# 1. Update the variant-allele-frequency / coefficient-of-variation figures.
$d.Content.Find.Execute(
    "VAFs of 5%, 10%-20%, 30%-40% and 50% are on average, 10.2%, 10.4%, 3.5% and 4.4%, respectively.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "VAFs of 2%, 5%-10%, 20%-40% and 50% are on average, 15.4%, 8.6%, 4.0% and 1.8%, respectively.",
    2)

# 2. Update the reported date (cached SAVEDATE field text).
$d.Content.Find.Execute(
    "1-Nov-2023", $true, $false, $false, $false, $false, $true, 1, $false,
    "16-Nov-2023", 2)
